$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A50").Value = 45236
$ws.Range("A50").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B50").Value = "20:24"
$ws.Range("C50").Value = 57.5
$ws.Range("D50").Value = "natura"

$ws.Range("A51").Value = 45237
$ws.Range("A51").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B51").Value = "20:26"
$ws.Range("C51").Value = 57.5
$ws.Range("D51").Value = "natura"
